$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric stay as text (matches source format)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.144.52'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '2.508.22'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '542.49'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("E6").Value = '  -2.88%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '2.533.16'
$ws.Range("E9").Value = '  +2.29%  '
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '5.53'
$ws.Range("E12").Value = '  +3.42%  '
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '2.954.15'
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("D15").Value = '23.58'
$ws.Range("E15").Value = '  -2.35%  '
$ws.Range("D16").Value = '59.069.89'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '2.522.78'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("D21").Value = '324.67'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +3.04%  '
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").Value = '61.82'
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("E25").Value = '  -5.05%  '
$ws.Range("E26").Value = '  +1.41%  '
$ws.Range("D27").Value = '0.994'
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").Value = '7.93'
$ws.Range("E28").Value = '  +2.38%  '
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").Value = '1.82'
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").Value = '1.21'
$ws.Range("E32").Value = '  -7.08%  '
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("E34").Value = '  +6.61%  '
$ws.Range("D35").Value = '158.19'
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("D36").Value = '18.66'
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("E38").Value = '  -6.65%  '
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("D40").Value = '36.92'
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("D41").Value = '297.98'
$ws.Range("E41").Value = '  -6.54%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").Value = '0.993'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '0.602'
$ws.Range("E45").Value = '  +2.89%  '
$ws.Range("D46").Value = '10.77'
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("D48").Value = '18.71'
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("D49").Value = '122.71'
$ws.Range("E49").Value = '  +0.78%  '
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("E51").Value = '  -0.91%  '

# Remove the temporary text number-format so styling matches the original (no s attribute)
$ws.Range("D2:D51").ClearFormats()

